# Update the "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the newly generated output numbers.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 183
    5  = 782
    6  = 232
    7  = 5753
    8  = 23
    10 = 96
    11 = 39
    12 = 27
    15 = 318
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
